$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B ("User ID") before the existing "Name" column,
# shifting Name, Reason, Amount, Account Number, Account Name, Bank Name
# one column to the right.
$ws.Columns("B:B").Insert()

# Update header row
$ws.Range("B1").Value = "User ID"

# Update data row 2 with the new/updated values
$ws.Range("A2").Value = "2025-05-23 13:17"
$ws.Range("B2").Value = "U07EVCPHEMP"
$ws.Range("C2").Value = "Gaius Omonale"
$ws.Range("E2").Value = 6000
$ws.Range("F2").Value = 7839920123
